$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. metsData sheet: insert two new rows (m2 and m8) that were previously
#    missing from the "mets" listing, shifting everything below them down.
# ---------------------------------------------------------------------------
$metsData = $wb.Worksheets.Item("metsData")

# Insert row for "m2" right after the "m3" row (old row 3 -> new row 4).
$metsData.Range("A4").EntireRow.Insert() | Out-Null
$metsData.Range("A4").Value2 = "m2"
$metsData.Range("B4").Value2 = 0.99
$metsData.Range("C4").Value2 = 1
$metsData.Range("D4").Value2 = 1.01

# Insert row for "m8" right after the "m7" row (now at row 8 after the shift
# above -> new row 9).
$metsData.Range("A9").EntireRow.Insert() | Out-Null
$metsData.Range("A9").Value2 = "m8"
$metsData.Range("B9").Value2 = 0.99
$metsData.Range("C9").Value2 = 1
$metsData.Range("D9").Value2 = 1.01

# ---------------------------------------------------------------------------
# 2. Column-width tweaks (cosmetic re-autofit deltas). The stored OOXML
#    "width" this runtime emits is ColumnWidth snapped to the nearest 1/6
#    plus 5/6, so we back-solve ColumnWidth = target - 5/6 to land as close
#    as possible to the target stored width.
# ---------------------------------------------------------------------------
$general = $wb.Worksheets.Item("general")
$general.Columns.Item(1).ColumnWidth = 72.863022941970266
$general.Columns.Item(2).ColumnWidth = 30.231443994601868

$measRates = $wb.Worksheets.Item("measRates")
$measRates.Columns.Item(1).ColumnWidth = 19.413630229419667

$kinetics1 = $wb.Worksheets.Item("kinetics1")
$kinetics1.Columns.Item(1).ColumnWidth = 15.020917678812365
$kinetics1.Columns.Item(2).ColumnWidth = 11.806342780026966
$kinetics1.Columns.Item(3).ColumnWidth = 11.915654520917666

$stoic = $wb.Worksheets.Item("stoic")
$stoic.Columns.Item(1).ColumnWidth = 12.020917678812467

$thermoRxns = $wb.Worksheets.Item("thermoRxns")
$thermoRxns.Columns.Item(2).ColumnWidth = 15.555330634277967
$thermoRxns.Columns.Item(3).ColumnWidth = 15.555330634277967

# ---------------------------------------------------------------------------
# 3. Selection / active-sheet bookkeeping.
#    Final state: "kinetics1" selection moves to H39 (without becoming the
#    active sheet), and "metsData" becomes the active sheet with A5 selected.
# ---------------------------------------------------------------------------
$kinetics1.Range("H39").Select() | Out-Null

$metsData.Activate()
$metsData.Range("A5").Select() | Out-Null
